$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New simulated values for columns A and B (rows 1-32)
$colA = @(
    -0.4092363652084714,
    -0.3621319999707353,
    -0.17617439202482998,
    -0.1634128991963575,
    -0.1567335908433991,
    -0.056605602655705134,
    -0.036565062324166675,
    -0.016502555770165017,
    -0.010472590710244312,
    -0.004449617009925078,
    0.000049514111616133505,
    0.006056492674837877,
    0.012114401044417633,
    0.004915768144987176,
    -0.02105126639014454,
    -0.015026831514303396,
    -0.009004082543484948,
    -0.03610850934001775,
    -0.027096515604499327,
    -0.018013308284090535,
    -0.00900425962852136,
    -0.09395379448802998,
    -0.08463905497186897,
    -0.0421276182854502,
    -0.09277408817599309,
    -0.08665967243362616,
    -0.08051750177986916,
    -0.07404821069419931,
    -0.061744873325194405,
    -0.04216844281024423,
    -0.027019175321632005,
    -0.00600063643732085
)
$colB = @(
    0.40824050861215255,
    0.35793857109562843,
    0.17541289901189572,
    0.16273359013198885,
    0.15537969010009878,
    0.056565061472094946,
    0.03650255491479015,
    0.016472589971765927,
    0.010449616270251205,
    0.0044504851615343455,
    -0.00005649341426661181,
    -0.00611440178195366,
    -0.01214565321123473,
    -0.004918491962696336,
    0.021026830776344596,
    0.015004081802822533,
    0.008999999233322065,
    0.036096514913928957,
    0.027013307587386848,
    0.018004258930901074,
    0.008999999301824158,
    0.09363905424696739,
    0.08412761726868112,
    0.04199999897778994,
    0.0926596717002468,
    0.08651750104371914,
    0.08004820994401385,
    0.07374487251466189,
    0.06161412786195619,
    0.04201917447397063,
    0.02700063553773724,
    0.005999999227757513
)

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Column widths adjusted to fit the new values (closest representable widths)
$ws.Range("A:A").ColumnWidth = 14.666666666666666
$ws.Range("B:B").ColumnWidth = 15.5

